$wb = $excel.ActiveWorkbook

$wsMeta  = $wb.Worksheets.Item("Metadata")
$wsCodes = $wb.Worksheets.Item("Include from Condition Inheri")

# --- Update the "Date" value (B8) ---
$wsMeta.Range("B8").Value = "2023-09-15T20:59:49+00:00"

# --- Update the "Jurisdiction" value (B12) ---
$wsMeta.Range("B12").Value = "Global (Whole world)"

# --- Re-apply wrap-text alignment on the existing header/body styles so the
#     stylesheet explicitly records applyAlignment="true" (the alignment
#     values themselves - vertical top / wrap text - are unchanged). Only
#     touch cells that already hold content so no new cells are created.
$wsMeta.Range("A1:B1").WrapText = $true
$wsMeta.Range("A2:B16").WrapText = $true

$wsCodes.Range("A1").WrapText = $true
$wsCodes.Range("A2").WrapText = $true
$wsCodes.Range("A3:B3").WrapText = $true
$wsCodes.Range("A4:B4").WrapText = $true
